# Generate Report for Handoff
# Update the "Latest Handoff Date(time)" values for the
# 5030e9b0-8513-4b90-ab57-3936d9dca066 row (row 6) on each sheet to reflect
# a newly generated handoff.

$wb = $excel.ActiveWorkbook

# Overview sheet - column D is "Latest Handoff Date"
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("D6").Value = "2016-28-12 12:28:57"

# zh-cn sheet - column E is "Latest Handoff Datetime"
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("E6").Value = "2016-03-12 12:28:54"

# de-de sheet - column E is "Latest Handoff Datetime"
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("E6").Value = "2016-03-12 12:28:57"
